$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Resolving-Mac" sending-cluster rows (old rows 14-16).
# This also drops the now-unused "Resolving-Mac" shared string on save.
$ws.Rows("14:16").Delete()

# Updated TPM-derived values for columns G:T (rows 2-13), one inner
# array per data row, in column order G,H,I,J,K,L,M,N,O,P,Q,R,S,T.
$newValues = @(
    @(0.3240676666666667,0.972203,0.04834775945726152,0.04834775945726153,1,0.3333333333333333,0.1574556666666667,0.472367,0.02266023449704293,0.02266023449704293,0.05102629050011111,0.459236614501,0.001095571566708171,0.001095571566708171),
    @(0.3240676666666667,0.972203,0.04834775945726152,0.04834775945726153,3,1,6.739756333333333,20.219269,0.9699521281096917,0.9699521281096917,2.184137108845222,19.657233979607,0.04689501217490628,0.04689501217490629),
    @(0.3240676666666667,0.972203,0.04834775945726152,0.04834775945726153,2,0.6666666666666666,0.05133333333333333,0.154,0.00738763739326543,0.00738763739326543,0.01663547355555556,0.149719262,0.0003571757156470675,0.0003571757156470676),
    @(4.372435,13.117305,0.6523249844605846,0.6523249844605848,1,0.3333333333333333,0.1574556666666667,0.472367,0.02266023449704293,0.02266023449704293,0.6884646678816667,6.196182010935,0.01478183711615673,0.01478183711615674),
    @(4.372435,13.117305,0.6523249844605846,0.6523249844605848,3,1,6.739756333333333,20.219269,0.9699521281096917,0.9699521281096917,29.46914648333833,265.222318350045,0.6327240068966656,0.6327240068966659),
    @(4.372435,13.117305,0.6523249844605846,0.6523249844605848,2,0.6666666666666666,0.05133333333333333,0.154,0.00738763739326543,0.00738763739326543,0.2244516633333334,2.02006497,0.004819140447762306,0.004819140447762307),
    @(0.036647,0.109941,0.005467377721001466,0.005467377721001468,1,0.3333333333333333,0.1574556666666667,0.472367,0.02266023449704293,0.02266023449704293,0.005770277816333333,0.051932500347,0.0001238920612418014,0.0001238920612418014),
    @(0.036647,0.109941,0.005467377721001466,0.005467377721001468,3,1,6.739756333333333,20.219269,0.9699521281096917,0.9699521281096917,0.2469918503476667,2.222926653129,0.005303094655664888,0.00530309465566489),
    @(0.036647,0.109941,0.005467377721001466,0.005467377721001468,2,0.6666666666666666,0.05133333333333333,0.154,0.00738763739326543,0.00738763739326543,0.001881212666666667,0.016930914,0.00004039100409477676,0.00004039100409477676),
    @(1.969698,5.909094,0.2938598783611522,0.2938598783611523,1,0.3333333333333333,0.1574556666666667,0.472367,0.02266023449704293,0.02266023449704293,0.310140111722,2.791261005498,0.006658933752936222,0.006658933752936222),
    @(1.969698,5.909094,0.2938598783611522,0.2938598783611523,3,1,6.739756333333333,20.219269,0.9699521281096917,0.9699521281096917,13.275284570254,119.477561132286,0.2850300143824547,0.2850300143824548),
    @(1.969698,5.909094,0.2938598783611522,0.2938598783611523,2,0.6666666666666666,0.05133333333333333,0.154,0.00738763739326543,0.00738763739326543,0.101111164,0.9100004759999999,0.002170930225761279,0.002170930225761279)
)

$cols = @(7,8,9,10,11,12,13,14,15,16,17,18,19,20)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $targetRow = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($targetRow, $cols[$j]).Value = $rowValues[$j]
    }
}
